# Update commonness index calculations (betas_hh.xlsx, sheet "HH models").
# Refreshed coefficient / bootstrap CI values for the household models,
# matching the 11/09 re-run of the commonness index pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HH models")

$updates = @{
    "D3"  = 0.38143458962440491
    "E3"  = 0.078820571303367615
    "D4"  = 0.39644512534141541
    "E4"  = 0.07859821617603302
    "D5"  = 0.4135589599609375
    "E5"  = 0.072071172297000885

    "H2"  = 0.16301143169403076
    "I2"  = 0.5299953818321228
    "H3"  = 0.22694626450538635
    "I3"  = 0.53592288494110107
    "H4"  = 0.24239261448383331
    "I4"  = 0.55049765110015869
    "H5"  = 0.27229946851730347
    "I5"  = 0.55481845140457153
    "H6"  = 0.85838252305984497
    "I6"  = 0.93988436460494995
    "H7"  = 0.86374819278717041
    "I7"  = 0.92985641956329346
    "H8"  = 0.86556857824325562
    "I8"  = 0.92760568857192993
    "H9"  = 0.86522281169891357
    "I9"  = 0.92518413066864014
    "H10" = 0.19332195818424225
    "I10" = 0.34843769669532776
    "H11" = 0.16804452240467072
    "I11" = 0.28901046514511108
    "H12" = 0.17637309432029724
    "I12" = 0.28824937343597412
    "H13" = 0.17728208005428314
    "I13" = 0.29531392455101013
    "H14" = 0.13654132187366486
    "I14" = 0.45524618029594421
    "H15" = 0.19173696637153625
    "I15" = 0.46316191554069519
    "H16" = 0.21224889159202576
    "I16" = 0.4762130081653595
    "H17" = 0.24249014258384705
    "I17" = 0.47845467925071716
    "H18" = -0.0065314383246004581
    "I18" = 0.13369239866733551
    "H19" = 0.015326528809964657
    "I19" = 0.13876241445541382
    "H20" = 0.0088002923876047134
    "I20" = 0.13419033586978912
    "H21" = 0.0091563351452350616
    "I21" = 0.13520373404026031
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
